# Add a new Time Log entry row for 11/10/2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 20) down to the new row (row 21)
$ws.Range("A20:C20").Copy() | Out-Null
$ws.Range("A21:C21").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new row's values
$ws.Range("A21").Value = (Get-Date -Year 2023 -Month 11 -Day 10 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("B21").Value = "Internship"
$ws.Range("C21").Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

# Update the active selection to the cell below the newly entered row, matching
# the natural behaviour of typing into C21 and pressing Enter.
$ws.Range("C22").Select() | Out-Null
